$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new date columns before column B (old B,C,D,E shift to E,F,G,H)
$ws.Columns("B").Insert()
$ws.Columns("B").Insert()
$ws.Columns("B").Insert()

# New header dates for the inserted columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns for every existing analyst row with the "UN" (unchanged) rank
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# New analyst rows added to the watch list
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
